$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the correlation matrix formulas (G5:J8) using CORREL, mirroring
# the header rows in G4:J4 / F5:F8 ("Column1".."Column4").
$dataCols = @('A', 'B', 'C', 'D')

for ($r = 0; $r -lt 4; $r++) {
    $row = 5 + $r
    $srcCol = $dataCols[$r]

    # Column G: single CORREL formula against A2:A11.
    $formulaG = '=CORREL(${0}$2:${0}$11, A2:A11)' -f $srcCol
    $ws.Range("G$row").Formula = $formulaG

    # Columns H:J: entered together so Excel stores them as one shared
    # formula (matching B2:B11 in the relative reference, auto-adjusted
    # for H/I/J -> B/C/D when copied across the row).
    $formulaH = '=CORREL(${0}$2:${0}$11, B2:B11)' -f $srcCol
    $ws.Range("H$row`:J$row").Formula = $formulaH
}

# Adjust the view to match: zoom level and active selection.
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("K11").Select()
